# Adaptations Data Model for Legal info on Resources
# Update deprecated functions
#
# - Rename property "hasCopyright"   -> "hasCopyrightResource"
# - Rename property "hasLicenseList" -> "hasLicenseResource"
# - Rename property "hasAuthorship"  -> "hasAuthorshipResource" (Archive sheet
#   still used the deprecated name; Material already used the Resource name)
# - Widen column A on the "Material" sheet so the longer property names fit
# - Refresh the remembered cell selection on both "Material" and "Archive"

$wb = $excel.ActiveWorkbook

# --- Rename the deprecated property names everywhere they are used -------
# xlWhole avoids "hasAuthorship" clobbering cells that already read
# "hasAuthorshipResource" (it would otherwise match as a substring and the
# replacement would be appended a second time).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("hasCopyright", "hasCopyrightResource", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
    $ws.Cells.Replace("hasLicenseList", "hasLicenseResource", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
    $ws.Cells.Replace("hasAuthorship", "hasAuthorshipResource", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
}

# --- Material sheet: widen column A and refresh the selection ------------
$material = $wb.Worksheets.Item("Material")
$material.Activate()
$material.Columns.Item(1).ColumnWidth = 33.5
$material.Range("A8").Select()

# --- Archive sheet: refresh the selection ---------------------------------
$archive = $wb.Worksheets.Item("Archive")
$archive.Activate()
$archive.Range("B14").Select()

# Leave "Material" as the active tab, matching the workbook's original state
$material.Activate()
